$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "Valor" column (B, numeric data rows 2-10) over to column C
# first, reading/writing cell-by-cell so the stored numeric values round-trip
# byte-for-byte (avoids Columns.Insert(), which reformats shifted values).
# Value2 is used on the read side since it returns the raw stored number/
# string without text-formatting loss.
for ($r = 2; $r -le 10; $r++) {
    $ws.Range("C$r").Value = $ws.Range("B$r").Value2
}

# Header row
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# New "Variável" column values (rows 2-10), overwriting the old column B data
# now that it has been copied to column C above.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Diferença 2022-2016"
}

# New "Colocação" (ranking) column values (rows 2-8 only)
$ws.Range("D2").Value = "1º"
$ws.Range("D3").Value = "2º"
$ws.Range("D4").Value = "3º"
$ws.Range("D5").Value = "4º"
$ws.Range("D6").Value = "5º"
$ws.Range("D7").Value = "6º"
$ws.Range("D8").Value = "27º"
